$d = $word.ActiveDocument

# --- Step 1: split the "Carrefour-Feuilles" run in the final paragraph,
#     wrapping "Feuilles" in spell-check proofErr markers.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$step1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0D649214" w14:textId="13A86F49" w:rsidR="00235BBD" w:rsidRDefault="00235BBD" w:rsidP="00235BBD"><w:r w:rsidRPr="00235BBD"><w:t>4.</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:tab/><w:t xml:space="preserve">Research the individual’s background, interests, and any relevant achievements or significant life events to provide context. Michel Patrick Mathurin is born in Haiti 1967 at Cabaret. Cabaret (Haitian Creole: Kabarè) is a commune in the Arcahaie Arrondissement, in the Ouest department of Haiti. Patrick favorite musical band in Tropicana, he </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>goes</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve"> in every single gala of them each time Tropicana is in Haiti for a tour. Patrick favorite beach is Moulin sur mer. Moulin Sur Mer Beach Resort is situated in the historical district of Montrouis, merely 1.9 km from Club Indigo and features 48 rooms with views of the sea. He loves to </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>read</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>actually</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve"> spend a lot of money in books. He's single, </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>with no</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve"> children. When he's in vacation </w:t></w:r><w:r w:rsidR="00793255"><w:t xml:space="preserve">is </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>in</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve"> the United States, Canada otherwise he's exploring Cap-Haitian, Jacmel, Port-Salut. Every morning, he does his 15 </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>minutes’ walk</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve"> before going to work. On </w:t></w:r><w:r w:rsidRPr="00235BBD"><w:lastRenderedPageBreak/><w:t xml:space="preserve">days off, he will not forgot his 3 </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>hours hiking</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve"> in  Fort Jacques, Petion-Ville Haiti. </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>Patrick</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve"> is the col</w:t></w:r><w:r w:rsidR="00793255"><w:t>l</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t>a</w:t></w:r><w:r w:rsidR="00793255"><w:t>bo</w:t></w:r><w:r w:rsidRPr="00235BBD"><w:t xml:space="preserve">rator of design of the </w:t></w:r><w:r w:rsidR="00793255" w:rsidRPr="00235BBD"><w:t>power plan</w:t></w:r><w:r><w:t xml:space="preserve"> substation of Carrefour-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Feuilles</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00793255"><w:t>.</w:t></w:r></w:p>
'@
$lastRange.InsertXML($step1Xml)

# --- Step 2: append a blank paragraph, then a new paragraph containing
#     the interview-prompt text (with proofErr markers around the
#     words Word's spell checker would flag).
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Give me 10 Interview question and answer to as Michel Patrick Mathurin, an electromechanical </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ngineer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, who live in Haiti, working at Electricity of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HAiti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. He like Tropicana Haiti, enjoy going to Moulin sur Mer for weekend at the beach, he likes to hike at </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Citadelle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LAferiere</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a cap </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>haitien</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Patrick was </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in Cabaret </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HAiti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in 1967.</w:t></w:r></w:p>
'@
$endRange.InsertXML($newParasXml)

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
